# Applies the "Added a few more slots" edit to the Double Diamond review:
#   1. Insert a new "Meta description" paragraph right after the Heading1 title.
#   2. Remove the duplicate bold "Play Double Diamond for Free - Classic Online Slot"
#      paragraph near the end of the document.
#   3. Replace the text of the final italic paragraph with the DALLE image prompt.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the duplicated bold title paragraph near the end of the document.
#    (Search starting right after the Heading1 title so we never match the
#    title itself, which has identical text.)
# ---------------------------------------------------------------------------
$titleEnd = $d.Paragraphs.Item(1).Range.End
$dupSearch = $d.Range($titleEnd, $d.Content.End)
$dupFound = $dupSearch.Find.Execute("Play Double Diamond for Free - Classic Online Slot", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($dupFound) {
    $dupPara = $dupSearch.Paragraphs.Item(1)
    $dupPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Replace the final italic paragraph's text with the DALLE image prompt
# ---------------------------------------------------------------------------
$oldBlurb = "Read our review and play Double Diamond for free at top online casinos. Enjoy the classic three-reel design and Wild symbol payouts."
$newPrompt = 'Prompt: DALLE, please create a cartoon-style feature image for the game "Double Diamond". The image should feature a happy Maya warrior with glasses. Make sure to incorporate the Double Diamond logo with the bright gems on either side and include the traditional symbols from old-fashioned mechanical slot machines, such as the numbers 7 and cherries. The image should be eye-catching and engaging to attract potential players.'

$blurbRange = $d.Content
$blurbFound = $blurbRange.Find.Execute($oldBlurb, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($blurbFound) {
    $blurbRange.Text = $newPrompt
}

# ---------------------------------------------------------------------------
# 3. Insert the "Meta description" paragraph after the Heading1 title
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$afterTitle = $titlePara.Range.Duplicate
$afterTitle.Collapse(0)
$afterTitle.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaRange = $d.Paragraphs.Item(2).Range
$metaRange.Text = "Meta description: Read our review and play Double Diamond for free at top online casinos. Enjoy the classic three-reel design and Wild symbol payouts."

# Make just the "Meta description" label bold (it is the first 16 characters
# of the paragraph we just created).
$labelStart = $metaPara.Range.Start
$labelRange = $d.Range($labelStart, $labelStart + 16)
$labelRange.Font.Bold = $true
